$wb = $excel.ActiveWorkbook

$wsRun = $wb.Worksheets.Item("RUNMANAGER")
$wsData = $wb.Worksheets.Item("DATA")
$wsBridge = $wb.Worksheets.Item("bridge_message_testdata")

# --- Data edits -----------------------------------------------------
# RUNMANAGER!B2: TC001_... -> loginLogoutTest
$wsRun.Range("B2").Value = "loginLogoutTest"

# DATA!B2 / DATA!C2: loginLogoutTest / QA (environment changed from DEV to QA)
$wsData.Range("B2").Value = "loginLogoutTest"
$wsData.Range("C2").Value = "QA"

# --- Column width edits ----------------------------------------------
# DATA column E gets much wider (used to show full QA description)
$wsData.Columns.Item(5).ColumnWidth = 36.6

# bridge_message_testdata gains a custom width on column B
$wsBridge.Columns.Item(2).ColumnWidth = 32.1

# --- Selections / active sheet ---------------------------------------
# RUNMANAGER: selection moves to B6
$wsRun.Range("B6").Select() | Out-Null

# bridge_message_testdata: selection becomes a single cell E25 (no longer tab-selected)
$wsBridge.Range("E25").Select() | Out-Null

# DATA becomes the active/visible sheet with selection C1
$wsData.Activate() | Out-Null
$wsData.Range("C1").Select() | Out-Null
